# "A lot of updates." - apply the set of changes to car.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B2/B4 ("Tim") become bold
$ws.Range("B2").Font.Bold = $true
$ws.Range("B4").Font.Bold = $true

# A6 ("X1") is retyped as text " 1" (leading space) using a Text number format
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = " 1"

# Page setup: paper size 9 (A4), portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to B4
[void]$ws.Range("B4").Select()
